$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Consolidate "The" + " " runs into a single "The " run.
$c1 = $tr.Characters(1, 4)
$c1.Text = "The "

# Consolidate "picture" + " " runs into a single "picture " run.
$c2 = $tr.Characters(5, 8)
$c2.Text = "picture "
